$d = $word.ActiveDocument

# 1) Simple text tweak: "Naming (below)." -> "Naming / Encodings (below)."
$d.Content.Find.Execute(
    "Naming (below).", $false, $false, $false, $false, $false,
    $true, 1, $false, "Naming / Encodings (below).", 2
) | Out-Null

# 2) Split the "Deep ML Embeddings..." paragraph into two list paragraphs:
#    - keep "...Behavior: regression. " in the original paragraph
#    - move "Naming: Auto Encoders." into a brand-new paragraph, which also
#      gains the extra trailing sentences from the commit.
#
# Using Find to locate the split point and InsertParagraphBefore (rather
# than a "^p" Find/Replace) keeps the existing run's rPr (the <w:rtl/>
# formatting) intact on both sides of the split.
$splitPoint = $d.Content
$found = $splitPoint.Find.Execute(
    "Naming: Auto Encoders.", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
$splitPoint.End = $splitPoint.Start
$splitPoint.InsertParagraphBefore()

# Append the rest of the new sentence onto the freshly created paragraph.
$tail = $d.Content
$tail.Find.Execute(
    "Naming: Auto Encoders.", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null
$tail.Collapse(0)
$tail.InsertAfter(" Semantic Hashing. Resources Mappings / Transforms Reified Maps / Tables. Keys / Values Resource Hashing / Resolution Functions: Contextual to Functional Environment State: Mappings Flows / Wrapped State.")
